$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# 1) Title paragraph: pStyle Title -> Heading1 (keep center alignment)
$pTitle = Get-ParaByText $d "AMIE and Usage Reporting v1"
$pTitle.Style = "Heading 1"
$pTitle.Alignment = 1

# 2) Subtitle paragraph -> drop Subtitle style + bookmark, retext, resize,
#    and add a new blank paragraph right after it (bold/sz24 paragraph mark).
$pSub = Get-ParaByText $d "Integration Roadmap Task"
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml2 = "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr/></w:pPr>" +
        "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/><w:rtl w:val=`"0`"/></w:rPr>" +
        "<w:t xml:space=`"preserve`">Infrastructure Integration Roadmap Task</w:t></w:r>" +
        "<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr></w:r></w:p>" +
        "<w:p $wNs><w:pPr><w:rPr><w:b w:val=`"1`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr></w:r></w:p>"
$pSub.Range.InsertXML($xml2)

# 3) Heading1 -> Heading2 for the four section headings that became subheadings.
foreach ($title in @("Summary", "Prerequisite tasks", "Support Information", "Detailed Instructions")) {
    $p = Get-ParaByText $d $title
    $p.Style = "Heading 2"
}

# 4) Merge the "Resource Providers..." paragraph with the following blank
#    paragraph, so the text run gains a trailing empty run and the blank
#    paragraph disappears.
$pSummaryBody = Get-ParaByText $d "Resource Providers will implement a client for the Account Management Information Exchange (AMIE) protocol, to receive allocations information from ACCESS Allocations and report necessary information back to ACCESS.  "
$nextPara = $pSummaryBody.Next()
$rngMerge = $d.Range($pSummaryBody.Range.Start, $nextPara.Range.End)
$xmlMerge = "<w:p $wNs><w:pPr><w:rPr/></w:pPr>" +
            "<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr>" +
            "<w:t xml:space=`"preserve`">Resource Providers will implement a client for the Account Management Information Exchange (AMIE) protocol, to receive allocations information from ACCESS Allocations and report necessary information back to ACCESS.  </w:t></w:r>" +
            "<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr></w:r></w:p>"
$rngMerge.InsertXML($xmlMerge)

# 5) Detailed Instructions heading: drop the trailing empty run.
$pDetailed = Get-ParaByText $d "Detailed Instructions"
$xmlDetailed = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/><w:rPr/></w:pPr>" +
               "<w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">Detailed Instructions</w:t></w:r></w:p>"
$pDetailed.Range.InsertXML($xmlDetailed)

# 6) Document Management heading: Heading1 -> Heading2.
$pDocMgmt = Get-ParaByText $d "Document Management"
$pDocMgmt.Style = "Heading 2"

Write-Output "done"
